# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Wed Jul 26 05:19:39 UTC 2023 with GitHub Actions".
#
# Rows 22/23 and 50/51 are full row-content swaps (coin name, link, price and
# volume all move to the adjacent row), every other touched row only updates
# its Price (D) and/or Volume(1h) (E) text.
#
# Price cells that look like plain numbers (e.g. "1.000", "0.5320") must stay
# literal text (leading/trailing zeros matter), so for those we briefly force
# a text NumberFormat, set the value, then restore the cell's original style
# so no new formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.227.65"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
$ws.Range("D3").Value = "1.856.78"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7054"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.97%  "

# Row 6
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.78"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08032"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +4.04%  "

# Row 9
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3020"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -0.33%  "

# Row 10
$ws.Range("E10").Value = "  +0.97%  "

# Row 11
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08172"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.79%  "

# Row 12
$ws.Range("D12").Value = "1.848.20"
$ws.Range("E12").Value = "  -0.25%  "

# Row 13
$ws.Range("E13").Value = "  -0.28%  "

# Row 14
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7033"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -2.87%  "

# Row 15
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.59"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").Value = "29.218.14"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007951"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +2.34%  "

# Row 18
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.792"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +0.97%  "

# Row 19
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.20"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.08%  "

# Row 20
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.44"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +0.99%  "

# Row 21
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.092.20"
$ws.Range("E22").Value = "  -0.25%  "

# Row 23
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.472"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -1.51%  "

# Row 25
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.82"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +1.43%  "

# Row 26
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.890"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.87%  "

# Row 27
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1430"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("E28").Value = "  +0.34%  "

# Row 29
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.917"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -3.03%  "

# Row 30
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.419"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +1.56%  "

# Row 31
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.475"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -0.70%  "

# Row 32
$ws.Range("E32").Value = "  -3.12%  "

# Row 33
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.012"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05186"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -0.67%  "

# Row 35
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.157"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -2.01%  "

# Row 36
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7134"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +1.82%  "

# Row 37
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -2.47%  "

# Row 38
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.650"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -0.19%  "

# Row 39
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01849"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -0.10%  "

# Row 40
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.721"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +1.66%  "

# Row 41
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9390"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +2.59%  "

# Row 42
$ws.Range("D42").Value = "1.128.86"
$ws.Range("E42").Value = "  +3.27%  "

# Row 43
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.935"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -1.08%  "

# Row 44
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4249"
$ws.Range("D44").Style = $origStyle

# Row 45
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.98"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -0.86%  "

# Row 46
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0000"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.82"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -0.01%  "

# Row 48
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5320"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -4.10%  "

# Row 49
$ws.Range("E49").Value = "  -0.55%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "1.991.18"
$ws.Range("E50").Value = "  -0.20%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.167"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.28%  "
